$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with refreshed market data ---

$ws.Range("D2").Value = "67.925.61"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "3.791.16"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.08"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.77"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "3.790.63"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.17"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "4.425.22"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "3.797.33"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "67.807.86"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.22"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.33"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.76"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.92"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.17"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "3.744.49"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.18"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.10"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "148.00"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "393.12"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.84"
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("D51").Value = "2.759.60"
$ws.Range("E51").Value = "  +2.66%  "

# --- Rows 26 and 27 swapped rank order: InternetComputer(DFINITY) moves up, Fetch.AI moves down ---

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.90"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  -3.15%  "
